$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 29 ("Click su item menu") content update --------------------------
# B29 (Espresso): swap the old navigationMenuItemView snippet for the updated
# overflow-menu "Change Settings" click snippet.
$b29 = @"
 openActionBarOverflowOrOptionsMenu(getInstrumentation().getTargetContext());
        ViewInteraction textView2 = onView(
                allOf(withId(android.R.id.title), withText("Change Settings"),
                        childAtPosition(
                                childAtPosition(
    withClassName(is("com.android.internal.view.menu.ListMenuItemView")),
                                        0),
                                0),
                        isDisplayed()));
        textView2.perform(click());
"@
$ws.Cells.Item(29, 2).Value = $b29

# C29 (Robolectric): replace the "//" placeholder with the real shadow call.
$ws.Cells.Item(29, 3).Value = "        shadowOf(activity).clickMenuItem(R.id.settings);"

# D29 (Robotium): replace the old clickOnMenuItem snippet with the new one
# and enable wrap text on the cell (new wrap-text style sharing font 4).
$d29 = @"
solo.sendKey(solo.MENU);
    // Click on Change Settings 
  solo.clickInList(3, 0);
"@
$ws.Cells.Item(29, 4).Value = $d29
$ws.Cells.Item(29, 4).WrapText = $true

# Row 29 grew taller to fit the new wrapped content.
$ws.Rows.Item(29).RowHeight = 195

# --- Column B got wider to fit the new text --------------------------------
$ws.Columns.Item(2).ColumnWidth = 72.5

# --- Selection moved to D37 -------------------------------------------------
$ws.Range("D37").Select()
